$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Name = "46131-0004"

$ws.Range("A205").Value = "Stand: 06.07.2020 / 16:14:31"

$ws.Range("C193").Value = 19137355
$ws.Range("D193").Value = 4941670397
$ws.Range("E193").Value = 3955618
$ws.Range("F193").Value = 1916175746
$ws.Range("G193").Value = 4775905
$ws.Range("H193").Value = 2071593829
$ws.Range("I193").Value = 1896470
$ws.Range("J193").Value = 1192333457
$ws.Range("K193").Value = 29765348
$ws.Range("L193").Value = 10121773429
